# Saldo.xlsx update
#
# Net effect of the source edit (derived from the OOXML diff):
#   - The row  004214592 / MERG / 75897.23        is removed entirely.
#   - The rows 004567880 / LUANA    / 16609.9
#              004265173 / JULIA    / 9000
#              004212438 / KENIA    / 4452.39
#              004376853 / ALBERTO  / 2401.39
#              004515341 / BRUNO    / 1104.61      are removed entirely.
#   - The row  005701765 / F / 51497.16 keeps its Conta/Nome but its Saldo
#     becomes 497.16, and the row itself moves down so it now sits right
#     after 005141215 / KARINA / 512.18 (i.e. after the block that used to
#     follow it) instead of right after 008013889 / CAROLINA / 85009.75.
#
# All other rows keep their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the six rows that disappear outright. Highest row number first
#    so earlier deletions don't shift the row numbers we still need.
#      row 3  -> 004214592 MERG     75897.23
#      row 6  -> 004567880 LUANA    16609.9
#      row 7  -> 004265173 JULIA    9000
#      row 8  -> 004212438 KENIA    4452.39
#      row 12 -> 004376853 ALBERTO  2401.39
#      row 13 -> 004515341 BRUNO    1104.61
$rowsToDelete = @(13, 12, 8, 7, 6, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).EntireRow.Delete() | Out-Null
}

# After those deletions the "F" row (005701765) now sits at row 3, and the
# "KARINA" row (005141215) now sits at row 36.

# 2) Grab the F row's current values before moving it.
$fConta = $ws.Cells.Item(3, 1).Value2
$fNome  = $ws.Cells.Item(3, 2).Value2

# 3) Remove the F row from its old spot (row 3). Everything below shifts up
#    by one, so KARINA (was row 36) is now row 35.
$ws.Rows.Item(3).EntireRow.Delete() | Out-Null

# 4) Insert a fresh row right after KARINA's new position (row 35), i.e. at
#    row 36, and populate it with F's data and the new Saldo value.
$ws.Rows.Item(36).EntireRow.Insert() | Out-Null

$fCell = $ws.Cells.Item(36, 1)
$fCell.NumberFormat = "@"
$fCell.Value2 = $fConta
$ws.Cells.Item(36, 2).Value2 = $fNome
$ws.Cells.Item(36, 3).Value2 = 497.16
